$d = $word.ActiveDocument

$replacements = @(
    @("143÷8=17, 7", "950÷2=475, 0"),
    @("808÷3=269, 1", "816÷5=163, 1"),
    @("342÷8=42, 6", "145÷2=72, 1"),
    @("308÷5=61, 3", "937÷2=468, 1"),
    @("140÷6=23, 2", "189÷2=94, 1"),
    @("990÷9=110, 0", "243÷6=40, 3"),
    @("784÷5=156, 4", "181÷7=25, 6"),
    @("662÷9=73, 5", "860÷8=107, 4"),
    @("947÷8=118, 3", "202÷4=50, 2"),
    @("214÷8=26, 6", "102÷6=17, 0"),
    @("236÷8=29, 4", "466÷4=116, 2"),
    @("508÷2=254, 0", "637÷8=79, 5"),
    @("899÷4=224, 3", "246÷8=30, 6"),
    @("482÷6=80, 2", "779÷8=97, 3"),
    @("437÷7=62, 3", "910÷9=101, 1"),
    @("291÷3=97, 0", "162÷4=40, 2"),
    @("299÷4=74, 3", "900÷6=150, 0"),
    @("667÷4=166, 3", "742÷4=185, 2"),
    @("292÷6=48, 4", "596÷9=66, 2"),
    @("203÷2=101, 1", "900÷7=128, 4"),
    @("882÷8=110, 2", "132÷4=33, 0"),
    @("731÷7=104, 3", "573÷4=143, 1"),
    @("589÷7=84, 1", "844÷4=211, 0"),
    @("201÷7=28, 5", "668÷6=111, 2"),
    @("957÷6=159, 3", "664÷5=132, 4")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
